$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room for the new rows: the worker-period table currently
#    occupies rows 16-39 (24 rows); the new table occupies rows
#    16-43 (28 rows). Insert 4 blank rows right after the last data
#    row (39) and before the blank gap / footer rows (44-45), which
#    pushes the footer (firm signature block) down to rows 48-49.
# ------------------------------------------------------------------
$ws.Rows.Item(40).Resize(4).Insert()

# Restore proper formatting on the newly inserted rows: rows 40-42
# should look like a normal data row (same formatting as row 38),
# and row 43 (new last data row) should carry the special "bottom of
# table" border formatting that used to belong to row 39.
$ws.Range("B38:J38").Copy()
$ws.Range("B40:J42").PasteSpecial(-4122)
$ws.Range("B39:J39").Copy()
$ws.Range("B43:J43").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Rewrite the worker-period detail table (rows 16-43) with the
#    refreshed data set: old/settled periods removed, new workers and
#    periods added, in the new order supplied by the refreshed
#    extract.
# ------------------------------------------------------------------
$rows = @(
  @{R=16; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2309"; F=88000;  G=2200000},
  @{R=17; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2308"; F=88000;  G=2200000},
  @{R=18; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2307"; F=88000;  G=2200000},
  @{R=19; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2306"; F=88000;  G=2200000},
  @{R=20; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2305"; F=88000;  G=2200000},
  @{R=21; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2304"; F=88000;  G=2200000},
  @{R=22; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2303"; F=88000;  G=2200000},
  @{R=23; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2302"; F=88000;  G=2200000},
  @{R=24; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2301"; F=88000;  G=2200000},
  @{R=25; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2212"; F=88000;  G=2200000},
  @{R=26; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2211"; F=88000;  G=2200000},
  @{R=27; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2210"; F=88000;  G=2200000},
  @{R=28; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2209"; F=88000;  G=2200000},
  @{R=29; B="CC"; C="79427140";   D="RAMON JOSE ARIZA RIOS";               E="2208"; F=88000;  G=2200000},
  @{R=30; B="CC"; C="73127846";   D="JOSE LUIS CARRILLO GRAU";             E="2302"; F=46400;  G=1160000},
  @{R=31; B="CC"; C="1047458914"; D="MERY CLARA MORELOS HENRIQUEZ";        E="2312"; F=120000; G=3000000},
  @{R=32; B="CC"; C="1047458914"; D="MERY CLARA MORELOS HENRIQUEZ";        E="2311"; F=120000; G=3000000},
  @{R=33; B="CC"; C="1047458914"; D="MERY CLARA MORELOS HENRIQUEZ";        E="2310"; F=120000; G=3000000},
  @{R=34; B="CC"; C="1050039005"; D="JUAN PABLO BUELVAS LEYVA";            E="1802"; F=31249;  G=781242},
  @{R=35; B="CC"; C="1052960119"; D="SHEYLA ANDREA PEREZ MIRANDA";         E="1709"; F=29509;  G=737717},
  @{R=36; B="CC"; C="1050969488"; D="ARLEY ALEXANDER MACIAS TRESPALACIOS"; E="2208"; F=40000;  G=1000000},
  @{R=37; B="CC"; C="1143401543"; D="ISAURA FILO ARBOLEDA";                E="1905"; F=6625;   G=828116},
  @{R=38; B="CC"; C="1052961171"; D="YESENIA ROCIO SCHMALBACH MORENO";     E="2505"; F=1898;   G=1423500},
  @{R=39; B="CC"; C="1052961171"; D="YESENIA ROCIO SCHMALBACH MORENO";     E="2503"; F=1898;   G=1423500},
  @{R=40; B="CC"; C="20373392";   D="LILIANA VANEGAS ORTEGA";              E="2207"; F=68658;  G=1980500},
  @{R=41; B="CC"; C="1052998110"; D="MARIAN ALEJANDRA BENAVIDEZ ACOSTA";   E="2208"; F=40000;  G=1000000},
  @{R=42; B="CC"; C="1007939014"; D="SHARON JULIANA GUERRERO RODRIGUEZ";   E="2304"; F=46400;  G=1160000},
  @{R=43; B="CC"; C="1007939014"; D="SHARON JULIANA GUERRERO RODRIGUEZ";   E="2303"; F=43307;  G=1160000}
)

foreach ($row in $rows) {
  $r = $row.R
  $ws.Range("B$r").Value = $row.B
  $ws.Range("C$r").Value = $row.C
  $ws.Range("D$r").Value = $row.D
  $ws.Range("E$r").Value = $row.E
  $ws.Range("F$r").Value = $row.F
  $ws.Range("G$r").Value = $row.G
}

# ------------------------------------------------------------------
# 3. Refresh the summary header figures to match the updated table:
#    total overdue amount, worker count and distinct-period count.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 1947944
$ws.Range("C13").Value = 11
$ws.Range("F13").Value = 23
